$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C255").Value = "Möchten Sie die Hub-Daten auch entfernen? Falls Sie später einmal auf Hub2 wechseln möchten, können die Daten übernommen werden. Allenfalls möchten Sie das Hub-Plugin nur deaktivieren?"
$ws.Range("D255").Value = "Do you want to remove the Hub data as well? If you want to switch to Hub2 later, the data can be transferred. At most, you just want to disable the Hub plugin?"

$ws.Range("C257").Value = "Entferne Hub-Daten"
$ws.Range("D257").Value = "Remove Hub data"

$ws.Range("C258").Value = "Hub-Daten behalten"
$ws.Range("D258").Value = "Keep Hub data"

$ws.Range("C259").Value = "Hub-Plugin nur deaktivieren"
$ws.Range("D259").Value = "Just deactivate Hub plugin"

$ws.Range("C260").Value = "Hub-Daten"
$ws.Range("D260").Value = "Hub data"

$ws.Range("C261").Value = "Die Hub-Daten wurden auch entfernt!"
$ws.Range("D261").Value = "The Hub data was also removed!"

$ws.Range("C262").Value = "Die Hub-Daten wurden behalten!"
$ws.Range("D262").Value = "The Hub data was kept!"

$ws.Range("B262").Select()
$excel.ActiveWindow.ScrollRow = 241
